$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.623.69'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').Value = '3.192.28'
$ws.Range('E3').Value = '  -5.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.06'
$ws.Range('E5').Value = '  -6.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.42'
$ws.Range('E6').Value = '  -9.00%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.191.55'
$ws.Range('E8').Value = '  -5.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.455'
$ws.Range('E9').Value = '  -5.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.44'
$ws.Range('E10').Value = '  -6.63%  '
$ws.Range('E11').Value = '  -7.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('E12').Value = '  -5.72%  '
$ws.Range('D13').Value = '3.747.32'
$ws.Range('E13').Value = '  -5.13%  '
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.68'
$ws.Range('E15').Value = '  -8.22%  '
$ws.Range('D16').Value = '3.199.43'
$ws.Range('E16').Value = '  -4.89%  '
$ws.Range('D17').Value = '58.674.51'
$ws.Range('E17').Value = '  -3.89%  '
$ws.Range('E18').Value = '  -8.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.89'
$ws.Range('E19').Value = '  -7.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.17'
$ws.Range('E20').Value = '  -9.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.10'
$ws.Range('E21').Value = '  -9.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '358.04'
$ws.Range('E22').Value = '  -4.64%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.65'
$ws.Range('E24').Value = '  -7.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.514'
$ws.Range('E25').Value = '  -8.54%  '
$ws.Range('D26').Value = '3.337.34'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.170'
$ws.Range('E27').Value = '  -3.40%  '
$ws.Range('D28').Value = '0.0₃0949'
$ws.Range('E28').Value = '  -12.23%  '
$ws.Range('E29').Value = '  +0.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.02'
$ws.Range('E30').Value = '  -5.35%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.91'
$ws.Range('E32').Value = '  -8.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.98'
$ws.Range('E33').Value = '  -9.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.57'
$ws.Range('E34').Value = '  -5.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.19'
$ws.Range('E35').Value = '  -8.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.34'
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.87'
$ws.Range('E37').Value = '  -9.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.30'
$ws.Range('E38').Value = '  -7.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.41'
$ws.Range('E39').Value = '  -9.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.73'
$ws.Range('E40').Value = '  -11.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0701'
$ws.Range('E41').Value = '  -7.26%  '
$ws.Range('D42').Value = '3.224.90'
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.82'
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.710'
$ws.Range('E44').Value = '  -6.72%  '
$ws.Range('E45').Value = '  -4.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.98'
$ws.Range('E46').Value = '  -7.34%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.48'
$ws.Range('E48').Value = '  -7.93%  '
$ws.Range('D49').Value = '2.283.78'
$ws.Range('E49').Value = '  -8.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.22'
$ws.Range('E50').Value = '  -6.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.43'
$ws.Range('E51').Value = '  -10.00%  '
